$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update "想去人数" (want-to-go count) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3075
$ws1.Range("F4").Value = 111
$ws1.Range("F5").Value = 6813
$ws1.Range("F6").Value = 1840
$ws1.Range("F7").Value = 53
$ws1.Range("F11").Value = 138
$ws1.Range("F12").Value = 155

# Sheet "演出" (Performances) - update "想去人数" value
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 7

# Sheet "全部类型" (All types) - update "想去人数" values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3075
$ws4.Range("F3").Value = 7
$ws4.Range("F5").Value = 111
$ws4.Range("F6").Value = 6813
$ws4.Range("F7").Value = 1840
$ws4.Range("F8").Value = 53
$ws4.Range("F12").Value = 138
$ws4.Range("F13").Value = 155
